$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix F90:F93: convert text dates "19/09/2025" into real date values
# matching the existing date format used elsewhere in column F (style from F2)
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F90:F93").PasteSpecial(-4122) | Out-Null
$ws.Range("F90:F93").Value = 45919

# --- Append new rows 94-96 (admin consuming seringa/algodao/gazes in Enfermagem)
$ws.Range("A94").Value = "admin"
$ws.Range("B94").Value = "seringa"
$ws.Range("C94").Value = 123
$ws.Range("D94").Value = "Enfermagem"
$ws.Range("E94").Value = "14:29:24"
$ws.Range("F94").Value = "20/09/2025"

$ws.Range("A95").Value = "admin"
$ws.Range("B95").Value = "algodão"
$ws.Range("C95").Value = 33
$ws.Range("D95").Value = "Enfermagem"
$ws.Range("E95").Value = "14:29:24"
$ws.Range("F95").Value = "20/09/2025"

$ws.Range("A96").Value = "admin"
$ws.Range("B96").Value = "gazes"
$ws.Range("C96").Value = 3
$ws.Range("D96").Value = "Enfermagem"
$ws.Range("E96").Value = "14:29:24"
$ws.Range("F96").Value = "20/09/2025"

Write-Host "Done"
